$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:E7 and G2:G7 (F column unchanged)
$values = @{
    2 = @{ B = 1.455362044514542;  C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697;  G = 3.754798637575387 }
    3 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697;  G = 5.586269137925634 }
    4 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697;  G = 6.189590430959694 }
    5 = @{ B = 1.455362044514542;  C = 1.655778082260271; D = 3.537761648806719;  E = 0.4942365360607697;  G = 7.143138311642302 }
    6 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697;  G = 5.586269137925634 }
    7 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 3.537761648806719;  E = 10.19245300693656;   G = 18.67282528286833 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("B$row").Value = $rowVals.B
    $ws.Range("C$row").Value = $rowVals.C
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("E$row").Value = $rowVals.E
    $ws.Range("G$row").Value = $rowVals.G
}
